$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; existing data (A:E) shifts to (B:F)
$ws.Range("A1").EntireColumn.Insert()

# Header for the newly inserted ID column - copy formatting from the
# neighboring header cell (bold, centered, bordered) then set its text
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial($xlPasteFormats)
$ws.Range("A1").Value = "ID"
$excel.CutCopyMode = 0

# Row labels for the newly inserted ID column
$ids = @("Hb 2", "Hb 3", "S 24", "S 28", "Hb 107", "Hb 66", "Hb 69", "Hb 95", "Hb 99", "Hb 92", "Hb 40", "Hb 41", "S 11", "Hb 57", "S 21", "S 22", "S 3", "S 4", "S 5", "Hb 74", "Hb 79", "Hb 32", "S 15", "S 16")

for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}
